# Update the embedded build timestamp across the workbook for the new release
# build: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")

$a2Range = $aboutSheet.Range("A2")
$a2Range.Value2 = $a2Range.Value2.Replace($oldStamp, $newStamp)

$a6Range = $aboutSheet.Range("A6")
$a6Range.Value2 = $a6Range.Value2.Replace($oldStamp, $newStamp)

$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")
$usedRange = $dataSheet.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)  # column S = build_version
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains($oldStamp)) {
        $cell.Value2 = $val.Replace($oldStamp, $newStamp)
    }
}
